# Weekly data refresh: a new week's price record is inserted at row 9
# (right after the fixed header block in rows 1-8), pushing all the
# existing historical rows down by one. This matches the commit message
# "Fruta / hortaliza, semanal" (weekly fruit/vegetable update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:56 down to 10:57, leaving a blank row 9 for the new record.
$ws.Rows(9).Insert()

# Populate the new record in row 9 with this week's data.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44547
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = "Poroto granado"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 44000
$ws.Range("L9").Value = 45000
$ws.Range("M9").Value = 44500
$ws.Range("N9").Value = "$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1780
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
